$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-parsed as numbers by Excel (matches original inlineStr text cells).
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated price (D) and volume (E) values row by row.
$ws.Range("D2").Value = "27.714.96"
$ws.Range("E2").Value = "  -0.55%  "
$ws.Range("D3").Value = "1.849.34"
$ws.Range("E3").Value = "  -1.15%  "
$ws.Range("D4").Value = "1.012"
$ws.Range("E4").Value = "  -2.96%  "
$ws.Range("D5").Value = "319.66"
$ws.Range("E5").Value = "  -1.81%  "
$ws.Range("D6").Value = "1.012"
$ws.Range("E6").Value = "  -2.60%  "
$ws.Range("D7").Value = "0.4330"
$ws.Range("E7").Value = "  -2.43%  "
$ws.Range("D8").Value = "0.3760"
$ws.Range("E8").Value = "  -1.30%  "
$ws.Range("D9").Value = "0.07369"
$ws.Range("E9").Value = "  -1.58%  "
$ws.Range("D10").Value = "0.8830"
$ws.Range("E10").Value = "  -0.54%  "
$ws.Range("D11").Value = "21.63"
$ws.Range("E11").Value = "  -1.03%  "
$ws.Range("D12").Value = "1.869.95"
$ws.Range("E12").Value = "  -0.32%  "
$ws.Range("D13").Value = "6.746"
$ws.Range("E13").Value = "  -0.51%  "
$ws.Range("D14").Value = "5.469"
$ws.Range("E14").Value = "  -1.94%  "
$ws.Range("D15").Value = "0.07131"
$ws.Range("E15").Value = "  -1.69%  "
$ws.Range("D16").Value = "87.94"
$ws.Range("E16").Value = "  +4.63%  "
$ws.Range("E17").Value = "  -2.76%  "
$ws.Range("D18").Value = "0.000009021"
$ws.Range("E18").Value = "  -1.82%  "
$ws.Range("D19").Value = "1.011"
$ws.Range("E19").Value = "  -2.60%  "
$ws.Range("D20").Value = "15.54"
$ws.Range("E20").Value = "  -0.45%  "
$ws.Range("D21").Value = "27.721.56"
$ws.Range("E21").Value = "  -0.60%  "
$ws.Range("D22").Value = "5.270"
$ws.Range("E22").Value = "  -1.26%  "
$ws.Range("D23").Value = "11.20"
$ws.Range("E23").Value = "  -1.75%  "
$ws.Range("D24").Value = "2.091.70"
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("D25").Value = "2.035"
$ws.Range("E25").Value = "  +2.90%  "
$ws.Range("D26").Value = "155.73"
$ws.Range("E26").Value = "  -2.10%  "
$ws.Range("D27").Value = "18.58"
$ws.Range("E27").Value = "  -1.88%  "
$ws.Range("D28").Value = "2.140"
$ws.Range("E28").Value = "  +6.79%  "
$ws.Range("D29").Value = "5.413"
$ws.Range("E29").Value = "  +1.20%  "
$ws.Range("D30").Value = "120.37"
$ws.Range("E30").Value = "  +2.02%  "
$ws.Range("D31").Value = "0.08957"
$ws.Range("E31").Value = "  -1.72%  "
$ws.Range("D32").Value = "1.237"
$ws.Range("E32").Value = "  +1.20%  "
$ws.Range("D33").Value = "0.7810"
$ws.Range("E33").Value = "  -0.29%  "
$ws.Range("D34").Value = "4.578"
$ws.Range("E34").Value = "  -0.29%  "
$ws.Range("D35").Value = "2.923"
$ws.Range("E35").Value = "  -6.03%  "
$ws.Range("D36").Value = "1.013"
$ws.Range("E36").Value = "  -2.90%  "
$ws.Range("D37").Value = "1.142"
$ws.Range("E37").Value = "  -2.07%  "
$ws.Range("D38").Value = "0.05346"
$ws.Range("E38").Value = "  -0.51%  "
$ws.Range("D39").Value = "0.01975"
$ws.Range("E39").Value = "  -1.36%  "
$ws.Range("D40").Value = "7.165"
$ws.Range("E40").Value = "  +3.41%  "
$ws.Range("D41").Value = "2.863"
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("D42").Value = "0.5183"
$ws.Range("E42").Value = "  -0.68%  "
$ws.Range("D43").Value = "0.1682"
$ws.Range("E43").Value = "  -0.95%  "
$ws.Range("D44").Value = "8.987"
$ws.Range("E44").Value = "  +3.17%  "
$ws.Range("D45").Value = "110.73"
$ws.Range("E45").Value = "  +0.54%  "
$ws.Range("D46").Value = "10.68"
$ws.Range("E46").Value = "  -0.88%  "
$ws.Range("D47").Value = "0.4744"
$ws.Range("E47").Value = "  +0.27%  "
$ws.Range("D48").Value = "1.709"
$ws.Range("E48").Value = "  -1.69%  "
$ws.Range("D49").Value = "0.06500"
$ws.Range("E49").Value = "  +0.48%  "
$ws.Range("D50").Value = "1.013"
$ws.Range("E50").Value = "  -2.86%  "
$ws.Range("D51").Value = "1.896"
$ws.Range("E51").Value = "  -0.39%  "
